# Update results in column AF ("average" under doctor_MA) for rows 4-13
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AF4").Value  = 0.767
$ws.Range("AF5").Value  = 1
$ws.Range("AF6").Value  = 0.868
$ws.Range("AF7").Value  = 0.9429999999999999
$ws.Range("AF8").Value  = 0.926
$ws.Range("AF9").Value  = 0.8
$ws.Range("AF10").Value = 1
$ws.Range("AF11").Value = 1
$ws.Range("AF12").Value = 1.2
$ws.Range("AF13").Value = 1.6
